$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Row 4: climate_change_factor_gnrl_hydropower_availability -> 1.63369506732604 for J4:AS4
$ws.Range("J4:AS4").Value = 1.63369506732604

# Row 5: elasticity_gnrl_rate_occupancy_to_gdppc -> -0.0317660546140297 for J5:AS5
$ws.Range("J5:AS5").Value = -0.0317660546140297
